$wb = $excel.ActiveWorkbook

$wsReactions = $wb.Worksheets.Item("Reactions")
$wsReactions.Activate()
$wsReactions.Range("H1").Value = "Flux bound units"
$wsReactions.Range("H1").Select()

$wsObjectives = $wb.Worksheets.Item("dFBA objectives")
$wsObjectives.Activate()
$wsObjectives.Range("F1:G1").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftToRight)
$wsObjectives.Range("F1").Value = "Reaction rate units"
$wsObjectives.Range("G1").Value = "Coefficient units"
$wsObjectives.Range("F1").Select()

Write-Output "done"
